# Weekly update: a new price observation is inserted as row 122
# (Femacal de La Calera / Ciboulette), pushing all the following
# rows down by one. The previously-last row (191) becomes row 192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 122; Excel shifts rows
# 122:191 down to 123:192 and the sheet dimension grows to R192.
$ws.Rows("122:122").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A122").Value = 3
$ws.Range("B122").Value = "Femacal de La Calera"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44488
$ws.Range("E122").Value = 5
$ws.Range("F122").Value = 100112039
$ws.Range("G122").Value = "Ciboulette"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 150
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 1500
$ws.Range("N122").Value = "$/docena de atados"
$ws.Range("O122").Value = "Provincia de Quillota"
$ws.Range("P122").Value = 500
$ws.Range("Q122").Value = 3
$ws.Range("R122").Value = "Hortaliza"
